# The commit inserts one new daily price record for Brócoli (Macroferia
# Regional de Talca) into the middle of the sheet, right after the existing
# row 345. That pushes every following record down by one row (old row 346
# becomes row 347, old row 464 becomes row 465, and the sheet grows from
# A1:R464 to A1:R465).
#
# Using EntireRow.Insert() reproduces exactly that behaviour: it shifts all
# rows below the insertion point down by one (carrying their values/styles
# with them) and leaves a blank row in their place, which we then populate
# with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 346; rows 346..464 (old) become 347..465 (new).
$ws.Rows.Item(346).Insert()

# Populate the newly inserted row 346 with the new record.
$ws.Cells.Item(346, 1).Value = 5
$ws.Cells.Item(346, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(346, 3).Value = "Maule"
$ws.Cells.Item(346, 4).Value = 44900
$ws.Cells.Item(346, 5).Value = 7
$ws.Cells.Item(346, 6).Value = 100112023
$ws.Cells.Item(346, 7).Value = "Brócoli"
$ws.Cells.Item(346, 8).Value = "Sin especificar"
$ws.Cells.Item(346, 9).Value = "Primera"
$ws.Cells.Item(346, 10).Value = 5000
$ws.Cells.Item(346, 11).Value = 500
$ws.Cells.Item(346, 12).Value = 500
$ws.Cells.Item(346, 13).Value = 500
$ws.Cells.Item(346, 14).Value = "$/unidad"
$ws.Cells.Item(346, 15).Value = "Región del Maule"
$ws.Cells.Item(346, 16).Value = 500
$ws.Cells.Item(346, 17).Value = 1
$ws.Cells.Item(346, 18).Value = "Hortaliza"

# Make sure the new row's date cell carries the same number format as the
# other date cells in column D (style index 2 in the original file).
$ws.Cells.Item(346, 4).NumberFormat = $ws.Cells.Item(345, 4).NumberFormat()
